# ---------------------------------------------------------------------------
# Applies the "render website, remove theme (not needed) from docs" edit:
#   1. Replace the single-cell wrapper table that holds the sampling-
#      distribution plot with a plain paragraph (style FirstParagraph)
#      that just contains the picture.
#   2. Add a new "Abstract Title" paragraph style and tighten the spacing
#      on the existing "Abstract" style.
#   3. Add a new "Footnote Block Text" paragraph style (based on
#      "Footnote Text").
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Unwrap the image out of its table and drop the (now empty) caption
#    paragraph that used to live alongside it.
# ---------------------------------------------------------------------------

# The paragraph right before the table is the last question of the
# "Repeated Samples" section ("... same distribution?").
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*distribution?*") {
        $anchor = $para
    }
}

# Insert a fresh paragraph right after it -- this lands cleanly outside
# (before) the table, unlike inserting at the table's own range boundary.
$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()

$pictureXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:drawing><wp:inline><wp:extent cx="4876800" cy="3657600"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="" title="" id="28" name="Picture"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="week-7-sampling-activity_files/figure-docx/slope-statisics-1.png" id="29" name="Picture"/><pic:cNvPicPr><a:picLocks noChangeArrowheads="1" noChangeAspect="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId27"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="4876800" cy="3657600"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525"><a:noFill/><a:headEnd/><a:tailEnd/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
"@
$newPara.Range.InsertXML($pictureXml)

# Now drop the original table (old picture run + caption paragraph). It is
# the single-cell (1 row x 1 column) table that used to wrap the picture;
# find it structurally rather than relying on a stale shape-range lookup.
$oldTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if (($t.Rows.Count -eq 1) -and ($t.Columns.Count -eq 1)) {
        $oldTable = $t
    }
}
$oldTable.Delete()

# ---------------------------------------------------------------------------
# 2. Styles: add "Abstract Title", tighten "Abstract" spacing.
# ---------------------------------------------------------------------------

$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------------
# 3. Styles: add "Footnote Block Text" (based on "Footnote Text").
# ---------------------------------------------------------------------------

$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "edit applied"
